$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on just the cells being updated so numeric-looking
# strings (e.g. "1.00", "27.729.54") are preserved exactly as text -
# matching the source data type - without touching untouched cells formatting.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.729.54'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.26%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.647.14'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.42%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.13'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.61%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.529'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -2.01%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.41'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +2.25%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.258'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.86%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.91%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0891'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.07%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.881.11'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.55%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.647.53'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.49%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.65%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.559'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.35%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.67'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.72%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.725.57'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.54%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '232.01'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.96%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.64'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +3.18%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.13%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.07'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +10.73%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -3.53%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '149.99'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.86%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.39%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -2.96%  '
$ws.Range('B28').NumberFormat = '@'
$ws.Range('B28').Value = 'BinanceUSD'
$ws.Range('C28').NumberFormat = '@'
$ws.Range('C28').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.25%  '
$ws.Range('B29').NumberFormat = '@'
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').NumberFormat = '@'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.64'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +1.15%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.86%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.28%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.30'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.28%  '
$ws.Range('B33').NumberFormat = '@'
$ws.Range('B33').Value = 'Maker'
$ws.Range('C33').NumberFormat = '@'
$ws.Range('C33').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.443.60'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +3.56%  '
$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.15'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +2.02%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +2.29%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.47%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.570'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +2.18%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.880'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.05%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +1.21%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.884'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +12.65%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.03'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.43%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.16%  '
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '67.11'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +4.88%  '
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.58'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +2.61%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.45'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.86%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.25'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +2.11%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.790.16'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.42%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.73'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +5.23%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +3.07%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '85.32'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.90%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.47%  '
